$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly price records were added to the "Cebollín" (Vega Modelo de
# Temuco) dataset. They belong chronologically before the existing row 300
# (which is why the author inserted them at that position), pushing the
# rest of the table (old rows 300-418) down by two rows to 302-420.
$ws.Rows("300:301").Insert()

# New record 1 (new row 300)
$newRow300 = @(
    10,
    "Vega Modelo de Temuco",
    "La Araucanía",
    44795,
    9,
    100112037,
    "Cebollín",
    "Sin especificar",
    "Primera",
    50,
    10000,
    10000,
    10000,
    "`$/docena de paquetes",
    "Provincia de Cautín",
    833,
    12,
    "Hortaliza"
)
for ($c = 1; $c -le 18; $c++) {
    $ws.Cells.Item(300, $c).Value = $newRow300[$c - 1]
}

# New record 2 (new row 301)
$newRow301 = @(
    10,
    "Vega Modelo de Temuco",
    "La Araucanía",
    44795,
    9,
    100112037,
    "Cebollín",
    "Sin especificar",
    "Primera",
    120,
    7000,
    7000,
    7000,
    "`$/docena de paquetes",
    "Región Metropolitana",
    583,
    12,
    "Hortaliza"
)
for ($c = 1; $c -le 18; $c++) {
    $ws.Cells.Item(301, $c).Value = $newRow301[$c - 1]
}
